$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 7554123858867932

# Row 3 - RandomForestRegressor
$ws.Range("B3").Value = 0.01282841001456107
$ws.Range("C3").Value = 0.01307048106583435
$ws.Range("D3").Value = 149428882609795.6

# Row 4 - rename model and update values
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.01355165461992512
$ws.Range("C4").Value = 0.01464076978721899
$ws.Range("D4").Value = 80091477075694.92

# Row 5 - rename model and update values
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 69492622709956.71
$ws.Range("C5").Value = 113134339237962.3
$ws.Range("D5").Value = 244414718392994.9
